$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "dogma1"
$ws.Range("A3:A7").EntireRow.Delete() | Out-Null
$ws.Range("A2").Select() | Out-Null
